# Add Plague Marines (and a Plague Champion) to the "Templar Models" sheet,
# per commit: "Added Plague Marines to the spreadsheet, made them the default enemy."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templar Models")

# Row 13: Plague Marine  (M WS BS T W A Ld Sv ... special1)
$ws.Range("A13").Value = "Plague Marine"
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 5
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 7
$ws.Range("J13").Value = 3
$ws.Range("L13").Value = 12

# Row 14: Plague Champion
$ws.Range("A14").Value = "Plague Champion"
$ws.Range("B14").Value = 5
$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 2
$ws.Range("I14").Value = 8
$ws.Range("J14").Value = 3
$ws.Range("L14").Value = 12

# Selection moves back to A4 on this sheet, as in the committed workbook.
$ws.Range("A4").Select()
